$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '33.801.40'
Set-TextValue $ws 'E2' '  +8.63%  '
Set-TextValue $ws 'D3' '1.775.27'
Set-TextValue $ws 'E3' '  +4.71%  '
Set-TextValue $ws 'E4' '  -0.12%  '
Set-TextValue $ws 'D5' '224.98'
Set-TextValue $ws 'E5' '  +1.83%  '
Set-TextValue $ws 'D6' '0.557'
Set-TextValue $ws 'E6' '  +4.44%  '
Set-TextValue $ws 'D7' '0.998'
Set-TextValue $ws 'E7' '  -0.11%  '
Set-TextValue $ws 'D8' '30.71'
Set-TextValue $ws 'E8' '  +3.37%  '
Set-TextValue $ws 'D9' '46.50'
Set-TextValue $ws 'E9' '  +2.63%  '
Set-TextValue $ws 'E10' '  +3.91%  '
Set-TextValue $ws 'D11' '0.0660'
Set-TextValue $ws 'E11' '  +3.24%  '
Set-TextValue $ws 'E12' '  +1.30%  '
Set-TextValue $ws 'D13' '2.025.42'
Set-TextValue $ws 'D14' '1.771.26'
Set-TextValue $ws 'E14' '  +3.28%  '
Set-TextValue $ws 'D15' '0.629'
Set-TextValue $ws 'E15' '  +2.88%  '
Set-TextValue $ws 'D16' '33.766.98'
Set-TextValue $ws 'E16' '  +8.67%  '
Set-TextValue $ws 'D17' '10.00'
Set-TextValue $ws 'E17' '  -2.91%  '
Set-TextValue $ws 'D18' '4.19'
Set-TextValue $ws 'E18' '  +1.47%  '
Set-TextValue $ws 'D19' '68.49'
Set-TextValue $ws 'E19' '  +2.17%  '
Set-TextValue $ws 'D20' '251.73'
Set-TextValue $ws 'E20' '  +1.38%  '
Set-TextValue $ws 'D21' '0.0₃0737'
Set-TextValue $ws 'E21' '  +2.29%  '
Set-TextValue $ws 'E22' '  -0.34%  '
Set-TextValue $ws 'D23' '10.27'
Set-TextValue $ws 'E23' '  +1.51%  '
Set-TextValue $ws 'E24' '  -2.28%  '
Set-TextValue $ws 'E25' '  -0.95%  '
Set-TextValue $ws 'D26' '158.88'
Set-TextValue $ws 'E26' '  +0.27%  '
Set-TextValue $ws 'D27' '16.51'
Set-TextValue $ws 'E27' '  +3.43%  '
Set-TextValue $ws 'D28' '0.115'
Set-TextValue $ws 'E28' '  +1.79%  '
Set-TextValue $ws 'D29' '6.94'
Set-TextValue $ws 'E29' '  +3.09%  '
Set-TextValue $ws 'D30' '0.998'
Set-TextValue $ws 'E30' '  -0.14%  '
Set-TextValue $ws 'D31' '3.80'
Set-TextValue $ws 'E31' '  +5.23%  '
Set-TextValue $ws 'D32' '0.0514'
Set-TextValue $ws 'E32' '  +2.50%  '
Set-TextValue $ws 'E33' '  +3.19%  '
Set-TextValue $ws 'D34' '3.55'
Set-TextValue $ws 'E34' '  +5.24%  '
Set-TextValue $ws 'D35' '1.82'
Set-TextValue $ws 'E35' '  +4.60%  '
Set-TextValue $ws 'D36' '1.483.10'
Set-TextValue $ws 'E36' '  -2.45%  '
Set-TextValue $ws 'E37' '  +3.09%  '
Set-TextValue $ws 'D38' '0.635'
Set-TextValue $ws 'E38' '  +2.12%  '
Set-TextValue $ws 'B39' 'VeChain'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D39' '0.0185'
Set-TextValue $ws 'E39' '  +2.84%  '
Set-TextValue $ws 'B40' 'Aave'
Set-TextValue $ws 'C40' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D40' '83.36'
Set-TextValue $ws 'E40' '  +0.42%  '
Set-TextValue $ws 'E41' '  +2.08%  '
Set-TextValue $ws 'E42' '  +0.74%  '
Set-TextValue $ws 'D43' '0.886'
Set-TextValue $ws 'E43' '  +4.28%  '
Set-TextValue $ws 'D44' '2.09'
Set-TextValue $ws 'E44' '  +2.46%  '
Set-TextValue $ws 'D45' '0.0515'
Set-TextValue $ws 'E45' '  +2.21%  '
Set-TextValue $ws 'E46' '  +3.59%  '
Set-TextValue $ws 'D47' '1.925.59'
Set-TextValue $ws 'E47' '  +5.49%  '
Set-TextValue $ws 'D48' '5.76'
Set-TextValue $ws 'E48' '  +3.25%  '
Set-TextValue $ws 'D49' '0.999'
Set-TextValue $ws 'D50' '11.74'
Set-TextValue $ws 'E50' '  +14.06%  '
Set-TextValue $ws 'D51' '50.64'
Set-TextValue $ws 'E51' '  -2.90%  '
